# Weekly refresh of "Fruta / hortaliza" data: the weekly rebuild re-sorts
# the existing price records, so several rows (2,4,6,7,9,10,11,12) end up
# swapping their Fecha/Volumen/Precio/Origen values with each other. Rows
# 1 (headers) and 3,5,8 (unaffected records) are left untouched. We write
# the new per-cell values directly so the result matches the target state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44644
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 2786
$ws.Range("O2").Value = "Provincia de Chacabuco"
$ws.Range("P2").Value = 464

# Row 4
$ws.Range("D4").Value = 44643
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 2800
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 2911
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 485

# Row 6
$ws.Range("D6").Value = 44672
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3500
$ws.Range("M6").Value = 3286
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 548

# Row 7
$ws.Range("D7").Value = 44659
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 2722
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 454

# Row 9
$ws.Range("D9").Value = 44637
$ws.Range("J9").Value = 170
$ws.Range("K9").Value = 2800
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 2906
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 484

# Row 10
$ws.Range("D10").Value = 44658
$ws.Range("J10").Value = 180
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2778
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 463

# Row 11
$ws.Range("D11").Value = 44631
$ws.Range("J11").Value = 110
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = 3273
$ws.Range("O11").Value = "Provincia de Chacabuco"
$ws.Range("P11").Value = 546

# Row 12
$ws.Range("D12").Value = 44630
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 2722
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 454
